# Reorder the worksheets so that "nobles" becomes the first sheet, then
# rename all four sheets to their (0-based) position index as a string,
# matching the "Full multi-move gem selection" restructuring.

$wb = $excel.ActiveWorkbook

# Move "nobles" to be the very first sheet in the workbook.
$nobles = $wb.Worksheets.Item("nobles")
$nobles.Move($wb.Worksheets.Item(1))

# Rename sheets in their new tab order: nobles, tier1, tier2, tier3 -> 0, 1, 2, 3
$wb.Worksheets.Item(1).Name = "0"
$wb.Worksheets.Item(2).Name = "1"
$wb.Worksheets.Item(3).Name = "2"
$wb.Worksheets.Item(4).Name = "3"
